$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Task# 105 ReadValue implementation
# Update the property override for test case 101 from TestCaseNumber=101 to Groups=Regression
$ws.Range("D2").Value = "Groups=Regression"

# Leave the active selection on D10, as left by the editor
$ws.Activate()
$ws.Range("D10").Select()
